$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (17) with the same shape/content pattern as existing rows:
# Category | Description | Date | Amount
# All existing rows store every value (including dates and numbers) as plain
# text shared strings, so force text formatting first - otherwise Excel would
# auto-convert "2023-03-07" into a real date and "300.0" into a number.
$row = $ws.Range("A17:D17")
$row.NumberFormat = "@"

$ws.Range("A17").Value = "Food"
$ws.Range("B17").Value = "food"
$ws.Range("C17").Value = "2023-03-07"
$ws.Range("D17").Value = "300.0"

# Drop the temporary text format again so the new cells end up using the
# workbook's default (unstyled) cell format, matching the rest of the sheet.
$row.ClearFormats()
